# Apply updated values to columns D, E, F, G, H, J for rows 2-25
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ "D" = 0.1058422808890214; "E" = 0.1302817993333925; "F" = 2.878478165862106; "G" = 2.368769996223534; "H" = 1.696868061312728; "J" = 0.2151472402880046 }
    3 = @{ "D" = 0.1057099897462948; "E" = 0.1293762481430178; "F" = 2.750918434283108; "G" = 2.222640816841476; "H" = 1.636692919714079; "J" = 0.2087731488325062 }
    4 = @{ "D" = 0.105631608910782; "E" = 0.1288919000434845; "F" = 2.674709327376604; "G" = 2.134686660620076; "H" = 1.600961959854402; "J" = 0.2050120596131251 }
    5 = @{ "D" = 0.1056003683076687; "E" = 0.1287124574646796; "F" = 2.64417696227224; "G" = 2.099280920985336; "H" = 1.586703403558658; "J" = 0.2035173224340951 }
    6 = @{ "D" = 0.1055952226750865; "E" = 0.1286837417798736; "F" = 2.639138474748648; "G" = 2.093427924020773; "H" = 1.584353915396832; "J" = 0.2032714010194852 }
    7 = @{ "D" = 0.1056311847743245; "E" = 0.1288894075202123; "F" = 2.674295448035082; "G" = 2.134207411306051; "H" = 1.600768445849013; "J" = 0.2049917480358658 }
    8 = @{ "D" = 0.1057960665953246; "E" = 0.1299546375540537; "F" = 2.834052162487808; "G" = 2.318011996819564; "H" = 1.67586474095458; "J" = 0.2129175040230962 }
    9 = @{ "D" = 0.1061427228874194; "E" = 0.1326164825695493; "F" = 3.164481609022914; "G" = 2.692914156986433; "H" = 1.832970800450028; "J" = 0.2296922894000204 }
    10 = @{ "D" = 0.1064127074359433; "E" = 0.1349279334951774; "F" = 3.418265018242863; "G" = 2.977790851650354; "H" = 1.954674957977716; "J" = 0.2427996710749483 }
    11 = @{ "D" = 0.1065390964415354; "E" = 0.1360582219595372; "F" = 3.536236937462036; "G" = 3.109577001089917; "H" = 2.011467405619044; "J" = 0.248939710659613 }
    12 = @{ "D" = 0.1065874905544675; "E" = 0.1364976832433129; "F" = 3.581282987014646; "G" = 3.159807628917747; "H" = 2.033183669984112; "J" = 0.2512908683670503 }
    13 = @{ "D" = 0.1065770439750011; "E" = 0.136402526482847; "F" = 3.571564784942154; "G" = 3.148974891615239; "H" = 2.028497258030598; "J" = 0.2507833367247514 }
    14 = @{ "D" = 0.1065430670494329; "E" = 0.1360941466472703; "F" = 3.539935378736089; "G" = 3.113702908760615; "H" = 2.013249775342672; "J" = 0.2491326164060297 }
    15 = @{ "D" = 0.1065223252557637; "E" = 0.1359067491117685; "F" = 3.520610273166028; "G" = 3.092140582485115; "H" = 2.003937769424056; "J" = 0.2481249124801792 }
    16 = @{ "D" = 0.1064045211602309; "E" = 0.1348556618024794; "F" = 3.410606892642249; "G" = 2.96922336625795; "H" = 1.950992589575776; "J" = 0.24240202150051 }
    17 = @{ "D" = 0.1063331807879706; "E" = 0.1342311160740586; "F" = 3.343776744307434; "G" = 2.894387060466499; "H" = 1.918881687400244; "J" = 0.2389370388623036 }
    18 = @{ "D" = 0.1062924826548031; "E" = 0.1338793036319004; "F" = 3.305574844941731; "G" = 2.851549089412686; "H" = 1.900546521073124; "J" = 0.2369607392022033 }
    19 = @{ "D" = 0.1062787599070294; "E" = 0.1337614551878588; "F" = 3.292680760349441; "G" = 2.837079934219787; "H" = 1.894361450487565; "J" = 0.2362944442791246 }
    20 = @{ "D" = 0.1063407402845478; "E" = 0.1342968323354974; "F" = 3.350866326720222; "G" = 2.902332128429407; "H" = 1.922286027038183; "J" = 0.2393041633429362 }
    21 = @{ "D" = 0.1065530322611821; "E" = 0.136184413718599; "F" = 3.549215510681393; "G" = 3.124054210283646; "H" = 2.017722587295964; "J" = 0.2496167613588085 }
    22 = @{ "D" = 0.1066948957063101; "E" = 0.1374848268051814; "F" = 3.681024872415435; "G" = 3.270868098814447; "H" = 2.081323391753415; "J" = 0.2565087896801259 }
    23 = @{ "D" = 0.1066188884859898; "E" = 0.1367846235793309; "F" = 3.610473322205934; "G" = 3.192332840114602; "H" = 2.047264520717533; "J" = 0.2528162789516699 }
    24 = @{ "D" = 0.106337321651834; "E" = 0.1342670994488628; "F" = 3.347660443095833; "G" = 2.898739583346071; "H" = 1.920746533548197; "J" = 0.2391381373017936 }
    25 = @{ "D" = 0.1060463370008105; "E" = 0.131834385398875; "F" = 2.589881758271019; "G" = 2.589881758271019; "H" = 1.789387422476466; "J" = 0.2250191059075632 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}

$wb.Save()
